# Apply cryptos list update (prices and volume %) per commit:
# "Updated cryptos list on Sun Mar 12 11:25:29 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    # Force the cell to be written back as text (matching the original
    # inlineStr/shared-string cell type) instead of letting Excel's
    # automatic type inference turn numeric-looking strings (e.g.
    # "20.553.59", "18.20", "0.9870") into numbers. ClearFormats()
    # afterwards removes the temporary "@" text number-format so the
    # cell's style index is left unchanged (matching the original,
    # unstyled data cells).
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue $ws 'D2' '20.553.59'
Set-TextValue $ws 'E2' '  +1.58%  '
Set-TextValue $ws 'D3' '1.473.54'
Set-TextValue $ws 'E3' '  +2.45%  '
Set-TextValue $ws 'D4' '1.005'
Set-TextValue $ws 'E4' '  -0.33%  '
Set-TextValue $ws 'D5' '0.9621'
Set-TextValue $ws 'E5' '  +6.06%  '
Set-TextValue $ws 'D6' '277.36'
Set-TextValue $ws 'E6' '  -0.03%  '
Set-TextValue $ws 'D7' '0.3622'
Set-TextValue $ws 'E7' '  -0.61%  '
Set-TextValue $ws 'D8' '0.3074'
Set-TextValue $ws 'E8' '  -1.03%  '
Set-TextValue $ws 'D9' '39.52'
Set-TextValue $ws 'E9' '  +1.39%  '
Set-TextValue $ws 'D10' '1.076'
Set-TextValue $ws 'E10' '  +5.65%  '
Set-TextValue $ws 'D11' '0.06656'
Set-TextValue $ws 'E11' '  +1.89%  '
Set-TextValue $ws 'E12' '  -0.12%  '
Set-TextValue $ws 'D13' '5.513'
Set-TextValue $ws 'E13' '  +2.58%  '
Set-TextValue $ws 'D14' '18.20'
Set-TextValue $ws 'E14' '  +3.60%  '
Set-TextValue $ws 'D15' '0.9618'
Set-TextValue $ws 'E15' '  +1.99%  '
Set-TextValue $ws 'D16' '6.167'
Set-TextValue $ws 'E16' '  +1.88%  '
Set-TextValue $ws 'D17' '0.00001026'
Set-TextValue $ws 'E17' '  +0.89%  '
Set-TextValue $ws 'D18' '1.473.38'
Set-TextValue $ws 'E18' '  +2.34%  '
Set-TextValue $ws 'D19' '0.05937'
Set-TextValue $ws 'E19' '  +5.07%  '
Set-TextValue $ws 'D20' '68.95'
Set-TextValue $ws 'E20' '  +1.92%  '
Set-TextValue $ws 'D21' '5.513'
Set-TextValue $ws 'E21' '  +2.32%  '
Set-TextValue $ws 'D22' '14.56'
Set-TextValue $ws 'E22' '  +1.16%  '
Set-TextValue $ws 'D23' '11.20'
Set-TextValue $ws 'E23' '  +3.96%  '
Set-TextValue $ws 'D24' '2.266'
Set-TextValue $ws 'E24' '  +1.36%  '
Set-TextValue $ws 'D25' '20.554.56'
Set-TextValue $ws 'D26' '143.15'
Set-TextValue $ws 'E26' '  +4.11%  '
Set-TextValue $ws 'D27' '2.128'
Set-TextValue $ws 'E27' '  -1.51%  '
Set-TextValue $ws 'D28' '17.16'
Set-TextValue $ws 'E28' '  +1.40%  '
Set-TextValue $ws 'E29' '  +2.62%  '
Set-TextValue $ws 'D30' '114.21'
Set-TextValue $ws 'E30' '  +3.93%  '
Set-TextValue $ws 'D31' '3.901'
Set-TextValue $ws 'E31' '  +0.21%  '
Set-TextValue $ws 'D32' '0.08014'
Set-TextValue $ws 'E32' '  +4.37%  '
Set-TextValue $ws 'D33' '4.948'
Set-TextValue $ws 'E33' '  +3.05%  '
Set-TextValue $ws 'D34' '0.8036'
Set-TextValue $ws 'E34' '  +0.41%  '
Set-TextValue $ws 'D35' '1.509'
Set-TextValue $ws 'E35' '  +5.01%  '
Set-TextValue $ws 'D36' '1.215'
Set-TextValue $ws 'E36' '  +6.23%  '
Set-TextValue $ws 'D37' '0.05777'
Set-TextValue $ws 'E37' '  -2.54%  '
Set-TextValue $ws 'D38' '4.733'
Set-TextValue $ws 'E38' '  +2.13%  '
Set-TextValue $ws 'D39' '0.02056'
Set-TextValue $ws 'E39' '  +3.59%  '
Set-TextValue $ws 'D40' '0.9624'
Set-TextValue $ws 'E40' '  +5.36%  '
Set-TextValue $ws 'D41' '10.43'
Set-TextValue $ws 'E41' '  +2.74%  '
Set-TextValue $ws 'E42' '  +1.97%  '
Set-TextValue $ws 'E43' '  +5.09%  '
Set-TextValue $ws 'D44' '0.5283'
Set-TextValue $ws 'E44' '  +1.07%  '
Set-TextValue $ws 'B45' 'EnergySwap'
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D45' '12.25'
Set-TextValue $ws 'E45' '  +2.07%  '
Set-TextValue $ws 'B46' 'PancakeSwap'
Set-TextValue $ws 'C46' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws 'D46' '3.521'
Set-TextValue $ws 'E46' '  +0.26%  '
Set-TextValue $ws 'D47' '118.98'
Set-TextValue $ws 'E47' '  +0.46%  '
Set-TextValue $ws 'D48' '0.5206'
Set-TextValue $ws 'E48' '  +1.39%  '
Set-TextValue $ws 'D49' '1.820'
Set-TextValue $ws 'E49' '  +3.70%  '
Set-TextValue $ws 'D50' '0.06469'
Set-TextValue $ws 'E50' '  +2.35%  '
Set-TextValue $ws 'D51' '0.9865'
Set-TextValue $ws 'E51' '  -0.21%  '
